$wb = $excel.ActiveWorkbook

# The "想去人数" (want-to-go count) figures were refreshed for both the
# "展览" sheet and the "全部类型" sheet (which mirrors the same rows).
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 2211
    $ws.Range("F3").Value = 1674
    $ws.Range("F4").Value = 327
    $ws.Range("F6").Value = 706
    $ws.Range("F8").Value = 5775
}
